$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:N25 (one inner array per row, columns B..N in order)
$data = @(
    @(0.301071449564887,0.04442704824161581,0,0.4160274979489884,2.358060693081583,0.8069844721067341,0.883284444334052,0.7074196399734234,0.06160035768159844,0.3213400753073756,0,0.410848662624403,1.75843662132662),
    @(0.2702915342972858,0.03903058097564838,0,0.4083150742256478,2.342876940250164,0.8085717031098625,0.8883650480969862,0.7123854761517876,0.06214069376210807,0.2869203827925446,0,0.391878410712664,1.777394830047317),
    @(0.2514541555025005,0.03572023839339522,0,0.4037670655084185,2.334807663926497,0.8100555676451648,0.8918672087951762,0.7157935636849366,0.0624895045906313,0.2658444808655958,0,0.3804200679449821,1.789625611394847),
    @(0.2437935776282529,0.03437202587360844,0,0.4019609195103087,2.33183460567561,0.8107881711458447,0.8933906119795836,0.7172726317261571,0.06263593587951366,0.2572707174848574,0,0.3757984925437299,1.794758158173848),
    @(0.2425225084711542,0.03414820318040768,0,0.4016638634114074,2.331359971492702,0.8109175408621638,0.8936493855584189,0.7175236792530519,0.06266050974849868,0.2558479547944472,0,0.3750339744168656,1.795619375541483),
    @(0.2513507778513713,0.0357020527782197,0,0.4037425159717785,2.33476629183258,0.8100649300636604,0.8918873642237202,0.7158131455705465,0.06249146204863809,0.2657287916562581,0,0.3803575460558193,1.789694229899205),
    @(0.2904459252202685,0.04256569727083104,0,0.4133293865643211,2.352565141102701,0.8074259518584199,0.884956831251074,0.7090573162431468,0.0617831300622429,0.3094602940846016,0,0.4042684780679693,1.764850854474215),
    @(0.3675904135866119,0.05605078817680464,0,0.4336149587925817,2.397420240683417,0.8062997337750062,0.8744015272511234,0.6986600720804503,0.0605291779554884,0.3956697446272415,0,0.4526573554894284,1.720818728315653),
    @(0.4245530513908875,0.06597596829398356,0,0.4494244904391991,2.436456044401567,0.8079526957581606,0.8684966170006305,0.6927624541932502,0.05969003097064629,0.4592798057854282,0,0.489121996974653,1.691325617192485),
    @(0.4505274028913959,0.07049554821148263,0,0.456813458335013,2.455538336786262,0.8092459777209058,0.8662120179613737,0.6904583101158508,0.05932606831744769,0.4882763550435811,0,0.5059091058963006,1.678529051768372),
    @(0.4603718588252264,0.07220767158291608,0,0.4596397814100541,2.462954938686309,0.8098137509377921,0.8654046406870748,0.6896403081341447,0.05919079775352021,0.4992650417525795,0,0.5122944988203812,1.673772541605373),
    @(0.4582513051511512,0.07183890597964648,0,0.4590298253986589,2.461349164943286,0.8096879971232624,0.8655759552585209,0.6898140532718031,0.0592198171070879,0.49689806618494,0,0.5109180261350019,1.674792968913936),
    @(0.4513371431926316,0.07063639247380138,0,0.4570454152140684,2.456144685793802,0.809291123723284,0.8661444371042961,0.6903899193934393,0.05931488834293797,0.4891802358694122,0,0.5064338661674128,1.678135941003887),
    @(0.447103120224682,0.06989990373199362,0,0.4558335877780806,2.452981609335581,0.8090581954988778,0.8665001696681429,0.6907497571526378,0.0593734548211069,0.4844539167183655,0,0.5036908927839008,1.680195237221193),
    @(0.4228567794415881,0.0656806949900357,0,0.4489455651631644,2.435235633381282,0.8078790862478797,0.8686540020214579,0.6929206619712005,0.05971417426365244,0.457385994471764,0,0.4880289138730589,1.692174388671097),
    @(0.4079980053409997,0.06309351645221284,0,0.444770429652813,2.424688384281239,0.807294527347139,0.8700781644329254,0.6943494879443719,0.05992774504347542,0.4407958532716236,0,0.4784716930841313,1.69968208841223),
    @(0.399457465617246,0.06160586952783831,0,0.4423875611330317,2.418746569784105,0.8070092466756904,0.8709351026973309,0.6952069528407634,0.0600522579162881,0.4312593099379001,0,0.4729933879639887,1.704058668769843),
    @(0.3965667982567425,0.06110225164970018,0,0.4415839523409986,2.416756186316903,0.8069213990872157,0.8712317389308311,0.6955033951734464,0.06009470314084409,0.4280313829113425,0,0.4711417541305138,1.705550521867458),
    @(0.4095791459241696,0.06336888154770293,0,0.4452129596490266,2.42579825216815,0.8073514809328799,0.8699226479564146,0.6941936976679273,0.05990483696822002,0.442561317826403,0,0.4794871355882435,1.698876841110897),
    @(0.4533677713621103,0.07098958202080041,0,0.4576275175677722,2.457668195466525,0.8094055756774736,0.8659758928226609,0.6902192931440041,0.05928689430425349,0.4914469262835723,0,0.5077502010072692,1.677151605321274),
    @(0.4820357431484013,0.07597397721376353,0,0.4659059633547429,2.47960775713527,0.8112029843327235,0.8637330789314319,0.6879396497107493,0.05889792041081643,0.5234448667007712,0,0.5263876719334775,1.663473306422752),
    @(0.4667306861827001,0.07331336151290202,0,0.4614725422816548,2.467796549716951,0.810201983314002,0.8648993084353691,0.6891272311835479,0.05910416123046325,0.5063626393242089,0,0.5164253737272446,1.670726016179605),
    @(0.408864306112406,0.06324438974796465,0,0.4450128374078588,2.425296101062685,0.8073255740206662,0.8699928380598578,0.6942640182737492,0.05991518832571474,0.4417631470748518,0,0.4790280033465706,1.699240705574642),
    @(0.3466702736020295,0.05239976291805704,0,0.4279681185829816,2.384218919261585,0.8061697004587103,0.8769320978285293,0.7011672590320011,0.06085396643927465,0.3722997801232566,0,0.4394064287587582,1.732228884512683)
)

$startRow = 2
$startCol = 2  # column B
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $c = $startCol + $j
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
